$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.230379746835443
$ws.Range("C2").Value = 0.4936708860759494
$ws.Range("J2").Value = 0.01518987341772152
$ws.Range("P2").Value = 0.1468354430379747
$ws.Range("S2").Value = 0.1139240506329114
$ws.Range("B3").Value = 0.009803921568627451
$ws.Range("C3").Value = 0.02941176470588235
$ws.Range("J3").Value = 0.03431372549019608
$ws.Range("P3").Value = 0.6813725490196079
$ws.Range("S3").Value = 0.2450980392156863
$ws.Range("J4").Value = 0.03636363636363636
$ws.Range("P4").Value = 0.6909090909090909
$ws.Range("S4").Value = 0.2727272727272727
$ws.Range("B6").Value = 0.04977375565610859
$ws.Range("D6").Value = 0.01357466063348416
$ws.Range("F6").Value = 0.02714932126696833
$ws.Range("J6").Value = 0.3212669683257919
$ws.Range("O6").Value = 0.02262443438914027
$ws.Range("Q6").Value = 0.1447963800904978
$ws.Range("R6").Value = 0.05882352941176471
$ws.Range("S6").Value = 0.3619909502262443
$ws.Range("B7").Value = 0.1310043668122271
$ws.Range("D7").Value = 0.03930131004366812
$ws.Range("F7").Value = 0.02620087336244541
$ws.Range("J7").Value = 0.1790393013100437
$ws.Range("Q7").Value = 0.1790393013100437
$ws.Range("R7").Value = 0.07860262008733625
$ws.Range("S7").Value = 0.3668122270742358
$ws.Range("B8").Value = 0.1007194244604317
$ws.Range("D8").Value = 0.01199040767386091
$ws.Range("F8").Value = 0.06235011990407674
$ws.Range("J8").Value = 0.1390887290167866
$ws.Range("O8").Value = 0.01199040767386091
$ws.Range("Q8").Value = 0.1702637889688249
$ws.Range("R8").Value = 0.07673860911270983
$ws.Range("S8").Value = 0.4268585131894485
$ws.Range("B9").Value = 0.06550218340611354
$ws.Range("D9").Value = 0.01310043668122271
$ws.Range("F9").Value = 0.06986899563318777
$ws.Range("J9").Value = 0.1179039301310044
$ws.Range("O9").Value = 0.01746724890829694
$ws.Range("Q9").Value = 0.2445414847161572
$ws.Range("R9").Value = 0.1179039301310044
$ws.Range("S9").Value = 0.3537117903930131
$ws.Range("B10").Value = 0.1398071625344353
$ws.Range("D10").Value = 0.02617079889807163
$ws.Range("F10").Value = 0.06129476584022039
$ws.Range("J10").Value = 0.1260330578512397
$ws.Range("O10").Value = 0.009641873278236915
$ws.Range("Q10").Value = 0.215564738292011
$ws.Range("R10").Value = 0.06955922865013774
$ws.Range("S10").Value = 0.3519283746556474
$ws.Range("G11").Value = 0.1111111111111111
$ws.Range("J11").Value = 0.1082621082621083
$ws.Range("K11").Value = 0.1680911680911681
$ws.Range("L11").Value = 0.5982905982905983
$ws.Range("S11").Value = 0.01424501424501425
$ws.Range("G12").Value = 0.7429906542056075
$ws.Range("J12").Value = 0.2289719626168224
$ws.Range("K12").Value = 0.004672897196261682
$ws.Range("L12").Value = 0.009345794392523364
$ws.Range("S12").Value = 0.01401869158878505
$ws.Range("G13").Value = 0.7115384615384616
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.03846153846153846
$ws.Range("J14").Value = 1
$ws.Range("F15").Value = 0.01578947368421053
$ws.Range("H15").Value = 0.2
$ws.Range("I15").Value = 0.07894736842105263
$ws.Range("J15").Value = 0.3947368421052632
$ws.Range("K15").Value = 0.04736842105263158
$ws.Range("M15").Value = 0.005263157894736842
$ws.Range("O15").Value = 0.03684210526315789
$ws.Range("S15").Value = 0.2210526315789474
$ws.Range("F16").Value = 0.02212389380530973
$ws.Range("H16").Value = 0.168141592920354
$ws.Range("I16").Value = 0.07079646017699115
$ws.Range("J16").Value = 0.3982300884955752
$ws.Range("K16").Value = 0.1194690265486726
$ws.Range("M16").Value = 0.01327433628318584
$ws.Range("N16").Value = 0.004424778761061947
$ws.Range("O16").Value = 0.03539823008849557
$ws.Range("S16").Value = 0.168141592920354
$ws.Range("F17").Value = 0.01181102362204724
$ws.Range("H17").Value = 0.1811023622047244
$ws.Range("I17").Value = 0.09055118110236221
$ws.Range("J17").Value = 0.4035433070866142
$ws.Range("K17").Value = 0.1141732283464567
$ws.Range("M17").Value = 0.02952755905511811
$ws.Range("N17").Value = 0.001968503937007874
$ws.Range("O17").Value = 0.05905511811023622
$ws.Range("S17").Value = 0.1082677165354331
$ws.Range("F18").Value = 0.02631578947368421
$ws.Range("H18").Value = 0.1526315789473684
$ws.Range("I18").Value = 0.131578947368421
$ws.Range("J18").Value = 0.4421052631578947
$ws.Range("K18").Value = 0.1
$ws.Range("M18").Value = 0.01052631578947368
$ws.Range("O18").Value = 0.04736842105263158
$ws.Range("S18").Value = 0.08947368421052632
$ws.Range("F19").Value = 0.01957831325301205
$ws.Range("H19").Value = 0.1694277108433735
$ws.Range("I19").Value = 0.0963855421686747
$ws.Range("J19").Value = 0.3930722891566265
$ws.Range("K19").Value = 0.1310240963855422
$ws.Range("M19").Value = 0.02259036144578313
$ws.Range("O19").Value = 0.0625
$ws.Range("S19").Value = 0.105421686746988
